{"js": "// Replace the two-digit-divided-by-one-digit problems in the practice\n// table with a new set of problems (text-only change; each original\n// value occurs exactly once in the document, so exact-text search +\n// replace is unambiguous).\nconst replacements = [\n  [\"57\u00f77=\", \"31\u00f79=\"],\n  [\"73\u00f74=\", \"22\u00f72=\"],\n  [\"28\u00f73=\", \"18\u00f73=\"],\n  [\"50\u00f79=\", \"27\u00f72=\"],\n  [\"96\u00f72=\", \"88\u00f73=\"],\n  [\"77\u00f79=\", \"76\u00f78=\"],\n  [\"60\u00f75=\", \"41\u00f79=\"],\n  [\"48\u00f75=\", \"68\u00f77=\"],\n  [\"60\u00f72=\", \"38\u00f79=\"],\n  [\"41\u00f76=\", \"16\u00f73=\"],\n  [\"12\u00f75=\", \"69\u00f78=\"],\n  [\"63\u00f72=\", \"58\u00f79=\"],\n  [\"70\u00f78=\", \"19\u00f74=\"],\n  [\"36\u00f79=\", \"38\u00f74=\"],\n  [\"22\u00f72=\", \"93\u00f75=\"],\n  [\"45\u00f75=\", \"30\u00f75=\"],\n  [\"13\u00f76=\", \"16\u00f76=\"],\n  [\"97\u00f79=\", \"47\u00f75=\"],\n  [\"27\u00f73=\", \"81\u00f78=\"],\n  [\"45\u00f72=\", \"42\u00f74=\"],\n  [\"86\u00f72=\", \"84\u00f78=\"],\n  [\"72\u00f78=\", \"34\u00f78=\"],\n  [\"11\u00f72=\", \"17\u00f78=\"],\n  [\"34\u00f72=\", \"55\u00f77=\"],\n  [\"12\u00f79=\", \"29\u00f77=\"],\n];\n\nconst body = context.document.body;\n\n// Phase 1: issue every search against the ORIGINAL document text and\n// load the hits. All `search()` calls below are resolved as of the\n// next `context.sync()`, before any text is mutated, so a replacement\n// value that happens to equal another row's original value (e.g.\n// \"22\u00f72=\") cannot be re-matched by a later rule in this same pass.\nconst pending = [];\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  pending.push({ results, newText });\n}\nawait context.sync();\n\n// Phase 2: now that all original-text locations are resolved, apply\n// every replacement.\nfor (const { results, newText } of pending) {\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the two-digit-divided-by-one-digit problems in the practice\n# table with a new set of problems (text-only change).\n#\n# Each original value occurs exactly once in the document, but some\n# original values equal some OTHER cell's new value (e.g. \"73\u00f74=\" ->\n# \"22\u00f72=\" while a different, later cell already holds \"22\u00f72=\" and must\n# become \"93\u00f75=\"). Doing this as one big Find/Replace-All per pair\n# (in sequence) would let a later rule re-match text that a previous\n# rule just wrote. To avoid that, we go cell-by-cell: each table cell\n# is visited exactly once and its replacement is looked up from its\n# own original text, so a newly written value is never re-examined.\n\n$replacements = @{\n    \"57\u00f77=\" = \"31\u00f79=\"\n    \"73\u00f74=\" = \"22\u00f72=\"\n    \"28\u00f73=\" = \"18\u00f73=\"\n    \"50\u00f79=\" = \"27\u00f72=\"\n    \"96\u00f72=\" = \"88\u00f73=\"\n    \"77\u00f79=\" = \"76\u00f78=\"\n    \"60\u00f75=\" = \"41\u00f79=\"\n    \"48\u00f75=\" = \"68\u00f77=\"\n    \"60\u00f72=\" = \"38\u00f79=\"\n    \"41\u00f76=\" = \"16\u00f73=\"\n    \"12\u00f75=\" = \"69\u00f78=\"\n    \"63\u00f72=\" = \"58\u00f79=\"\n    \"70\u00f78=\" = \"19\u00f74=\"\n    \"36\u00f79=\" = \"38\u00f74=\"\n    \"22\u00f72=\" = \"93\u00f75=\"\n    \"45\u00f75=\" = \"30\u00f75=\"\n    \"13\u00f76=\" = \"16\u00f76=\"\n    \"97\u00f79=\" = \"47\u00f75=\"\n    \"27\u00f73=\" = \"81\u00f78=\"\n    \"45\u00f72=\" = \"42\u00f74=\"\n    \"86\u00f72=\" = \"84\u00f78=\"\n    \"72\u00f78=\" = \"34\u00f78=\"\n    \"11\u00f72=\" = \"17\u00f78=\"\n    \"34\u00f72=\" = \"55\u00f77=\"\n    \"12\u00f79=\" = \"29\u00f77=\"\n}\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nforeach ($row in $table.Rows) {\n    foreach ($cell in $row.Cells) {\n        # Cell.Range.Text includes the trailing cell-mark (CR + BEL);\n        # strip it so we can match against the plain problem text.\n        $cellText = $cell.Range.Text\n        $plainText = $cellText.TrimEnd([char]13, [char]7)\n\n        if ($replacements.ContainsKey($plainText)) {\n            $cell.Range.Text = $replacements[$plainText]\n        }\n    }\n}\n"}
